# Update "want to go" counts (column F) on the "展览" and "全部类型" sheets.
# Each listed event's count is bumped up by a small amount as captured in the
# source diff (a generated-output refresh).

$wb = $excel.ActiveWorkbook

$updates = @(
    @{F=2;  Old=388;   New=389},
    @{F=3;  Old=681;   New=682},
    @{F=4;  Old=114;   New=115},
    @{F=5;  Old=2111;  New=2113},
    @{F=6;  Old=8;     New=9},
    @{F=7;  Old=10986; New=10993},
    @{F=8;  Old=187;   New=188},
    @{F=9;  Old=166;   New=167},
    @{F=10; Old=298;   New=300},
    @{F=11; Old=211;   New=212},
    @{F=12; Old=10841; New=10855},
    @{F=16; Old=755;   New=756},
    @{F=17; Old=5438;  New=5443},
    @{F=18; Old=84;    New=85},
    @{F=19; Old=3406;  New=3408}
)

$ws1 = $wb.Worksheets.Item("展览")
foreach ($u in $updates) {
    $ws1.Range("F" + $u.F).Value = $u.New
}

$updates4 = @(
    @{F=2;  New=389},
    @{F=3;  New=682},
    @{F=5;  New=115},
    @{F=6;  New=2113},
    @{F=8;  New=9},
    @{F=10; New=10993},
    @{F=11; New=188},
    @{F=12; New=167},
    @{F=13; New=300},
    @{F=14; New=212},
    @{F=15; New=10855},
    @{F=19; New=756},
    @{F=20; New=5443},
    @{F=21; New=85},
    @{F=22; New=3408}
)

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($u in $updates4) {
    $ws4.Range("F" + $u.F).Value = $u.New
}
